# Update country data file: "Data" sheet -> "Summary" sheet, with expanded
# MSME indicator content (Source Type line, a second "Value added" table,
# and a Sector Distribution Details / source block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet -------------------------------------------------
$ws.Name = "Summary"

# --- Wipe out the old rows 5-13 content/formatting so we can rebuild it
#     at the new row numbers without leftover formatting artifacts.
$ws.Range("A5:D13").Clear()

# --- Row 7: new "Source Type" sub-header (bold + underline) -----------
$ws.Range("A7").Value = "Source Type: Statistical Institution (Most Widely Used)"
$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").Font.Underline = $true

# --- Row 9: column headers (bold) --------------------------------------
$ws.Range("B9").Value = "Micro"
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Value = "SMEs"
$ws.Range("C9").Font.Bold = $true
$ws.Range("D9").Value = "MSMEs"
$ws.Range("D9").Font.Bold = $true

# --- Row 10: Enterprises (absolute #) ----------------------------------
$ws.Range("A10").Value = "Enterprises (absolute #)"
$ws.Range("A10").Font.Bold = $true
$ws.Range("D10").Value = "'253080"
$ws.Range("D10").Style = "Normal"

# --- Row 11: Enterprises density (per 1000 people) ---------------------
$ws.Range("A11").Value = "Enterprises density (per 1000 people)"
$ws.Range("A11").Font.Bold = $true
$ws.Range("D11").Value = "'8.5"
$ws.Range("D11").Style = "Normal"

# --- Row 12: Employment (absolute #) ------------------------------------
$ws.Range("A12").Value = "Employment (absolute #)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("D12").Value = "'9519600"
$ws.Range("D12").Style = "Normal"

# --- Row 13: source note (italic) ---------------------------------------
$ws.Range("A13").Value = "Source: SCRUS, 2013"
$ws.Range("A13").Font.Italic = $true

# --- Row 15: second table's column headers (bold) -----------------------
$ws.Range("B15").Value = "Micro"
$ws.Range("B15").Font.Bold = $true
$ws.Range("C15").Value = "SMEs"
$ws.Range("C15").Font.Bold = $true
$ws.Range("D15").Value = "MSMEs"
$ws.Range("D15").Font.Bold = $true

# --- Row 16: Value added to the economy ---------------------------------
$ws.Range("A16").Value = "Value added to the economy (% of total)"
$ws.Range("A16").Font.Bold = $true
$ws.Range("D16").Value = "'51"
$ws.Range("D16").Style = "Normal"

# --- Row 17: source note (italic) ----------------------------------------
$ws.Range("A17").Value = "Source: SCRUS, 2013"
$ws.Range("A17").Font.Italic = $true

# --- Row 22: Sector Distribution Details (bold) ---------------------------
$ws.Range("A22").Value = "Sector Distribution Details"
$ws.Range("A22").Font.Bold = $true

# --- Row 25: SCRUS (bold) --------------------------------------------------
$ws.Range("A25").Value = "SCRUS"
$ws.Range("A25").Font.Bold = $true

# --- Row 26: long source citation text (italic) ----------------------------
$ws.Range("A26").Value = "The State Committee of the Republic of Uzbekistan on Statistics (SCRUS), ""По состоянию на 1 октября 2012 года количество зарегистрированных субъектов малого бизнеса (без фермерских хозяйств) составило 249 915"", 2012. Available at http://www.stat.uz/press/1/5154/?sphrase_id=108660`nThe State Committee of the Republic of Uzbekistan on Statistics (SCRUS), ""Об итогах социально-экономического развития Республики Узбекистан за I квартал 2014 года"", 2014. Available at http://www.stat.uz/press/1/8359/?sphrase_id=108660"
$ws.Range("A26").Font.Italic = $true
